{"js": "// Resume update: mark \"Convolutional Neural Networks\" as completed\n// (move it from the \"In Progress\" list to the end of the \"Completed\" list),\n// and drop the stray \"_GoBack\" bookmark left over from the last edit.\n\nconst body = context.document.body;\n\n// 1. Remove the leftover \"_GoBack\" bookmark (Word auto-bookmark of the\n//    last edit location; not meaningful content, just edit residue).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2. Append \"; Convolutional Neural Networks\" to the end of the\n//    \"Completed\" courses line.\nconst completedResults = body.search(\"Structuring Machine Learning Projects\", { matchCase: true });\ncompletedResults.load(\"text\");\nawait context.sync();\n\nif (completedResults.items.length > 0) {\n  const completedRange = completedResults.items[0];\n  completedRange.insertText(\"; Convolutional Neural Networks\", Word.InsertLocation.after);\n}\n\n// 3. Remove \"Convolutional Neural Networks; \" from the start of the\n//    \"In Progress\" courses line, leaving just \"Sequence Models\".\nconst inProgressResults = body.search(\"Convolutional Neural Networks; Sequence Models\", { matchCase: true });\ninProgressResults.load(\"text\");\nawait context.sync();\n\nif (inProgressResults.items.length > 0) {\n  inProgressResults.items[0].insertText(\"Sequence Models\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Resume update: mark \"Convolutional Neural Networks\" as completed\n# (move it from the \"In Progress\" list to the end of the \"Completed\" list),\n# and drop the stray \"_GoBack\" bookmark left over from the last edit.\n\n$d = $word.ActiveDocument\n\n# 1. Remove the leftover \"_GoBack\" bookmark (Word auto-bookmark of the\n#    last edit location; not meaningful content, just edit residue).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Append \"; Convolutional Neural Networks\" to the end of the\n#    \"Completed\" courses line.\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Text = \"Structuring Machine Learning Projects\"\n$range.Find.Replacement.Text = \"Structuring Machine Learning Projects; Convolutional Neural Networks\"\n$range.Find.Forward = $true\n$range.Find.Wrap = 0\n[void]$range.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n# 3. Remove \"Convolutional Neural Networks; \" from the start of the\n#    \"In Progress\" courses line, leaving just \"Sequence Models\".\n$range2 = $d.Content\n$range2.Find.ClearFormatting()\n$range2.Find.Text = \"Convolutional Neural Networks; Sequence Models\"\n$range2.Find.Replacement.Text = \"Sequence Models\"\n$range2.Find.Forward = $true\n$range2.Find.Wrap = 0\n[void]$range2.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
